$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (62) into the new rows.
# Column A uses style s=1 (bold, bordered, centered) and column E uses
# style s=2 (custom date/time number format). Using Copy + PasteSpecial
# (xlPasteFormats) reuses the existing style indices instead of creating
# new duplicate styles.
$xlPasteFormats = -4122

# Row 63 (Indice=62)
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial($xlPasteFormats)
$ws.Range("E62").Copy()
$ws.Range("E63").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(63, 1).Value = 62
$ws.Cells.Item(63, 2).Value = "spain"
$ws.Cells.Item(63, 3).Value = "laliga"
$ws.Cells.Item(63, 4).Value = "2023-2024"
$ws.Cells.Item(63, 5).Value = 45196.79166666666
$ws.Cells.Item(63, 6).Value = "Ath Bilbao"
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = "Getafe"
$ws.Cells.Item(63, 9).Value = 2
$ws.Cells.Item(63, 10).Value = 1.71
$ws.Cells.Item(63, 11).Value = "17/09/2023 09:02"
$ws.Cells.Item(63, 12).Value = 1.53
$ws.Cells.Item(63, 13).Value = "27/09/2023 18:31"
$ws.Cells.Item(63, 14).Value = 3.42
$ws.Cells.Item(63, 15).Value = "17/09/2023 09:02"
$ws.Cells.Item(63, 16).Value = 4.06
$ws.Cells.Item(63, 17).Value = "27/09/2023 18:49"
$ws.Cells.Item(63, 18).Value = 5.44
$ws.Cells.Item(63, 19).Value = "17/09/2023 09:02"
$ws.Cells.Item(63, 20).Value = 7.73
$ws.Cells.Item(63, 21).Value = "27/09/2023 18:49"
$ws.Cells.Item(63, 22).Value = "https://www.betexplorer.com/football/spain/laliga/ath-bilbao-getafe/zgsFCYIT/"

# Row 64 (Indice=63)
$ws.Range("A62").Copy()
$ws.Range("A64").PasteSpecial($xlPasteFormats)
$ws.Range("E62").Copy()
$ws.Range("E64").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(64, 1).Value = 63
$ws.Cells.Item(64, 2).Value = "spain"
$ws.Cells.Item(64, 3).Value = "laliga"
$ws.Cells.Item(64, 4).Value = "2023-2024"
$ws.Cells.Item(64, 5).Value = 45196.79166666666
$ws.Cells.Item(64, 6).Value = "Real Madrid"
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 8).Value = "Las Palmas"
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 1.2
$ws.Cells.Item(64, 11).Value = "23/09/2023 09:28"
$ws.Cells.Item(64, 12).Value = 1.18
$ws.Cells.Item(64, 13).Value = "27/09/2023 18:29"
$ws.Cells.Item(64, 14).Value = 6.76
$ws.Cells.Item(64, 15).Value = "23/09/2023 09:28"
$ws.Cells.Item(64, 16).Value = 8
$ws.Cells.Item(64, 17).Value = "27/09/2023 18:29"
$ws.Cells.Item(64, 18).Value = 11.3
$ws.Cells.Item(64, 19).Value = "23/09/2023 09:28"
$ws.Cells.Item(64, 20).Value = 16.5
$ws.Cells.Item(64, 21).Value = "27/09/2023 18:29"
$ws.Cells.Item(64, 22).Value = "https://www.betexplorer.com/football/spain/laliga/real-madrid-las-palmas/GQHmRXXM/"

# Row 65 (Indice=64)
$ws.Range("A62").Copy()
$ws.Range("A65").PasteSpecial($xlPasteFormats)
$ws.Range("E62").Copy()
$ws.Range("E65").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(65, 1).Value = 64
$ws.Cells.Item(65, 2).Value = "spain"
$ws.Cells.Item(65, 3).Value = "laliga"
$ws.Cells.Item(65, 4).Value = "2023-2024"
$ws.Cells.Item(65, 5).Value = 45196.79166666666
$ws.Cells.Item(65, 6).Value = "Villarreal"
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = "Girona"
$ws.Cells.Item(65, 9).Value = 2
$ws.Cells.Item(65, 10).Value = 1.71
$ws.Cells.Item(65, 11).Value = "17/09/2023 09:02"
$ws.Cells.Item(65, 12).Value = 2.17
$ws.Cells.Item(65, 13).Value = "27/09/2023 18:51"
$ws.Cells.Item(65, 14).Value = 4.14
$ws.Cells.Item(65, 15).Value = "17/09/2023 09:02"
$ws.Cells.Item(65, 16).Value = 3.72
$ws.Cells.Item(65, 17).Value = "27/09/2023 18:51"
$ws.Cells.Item(65, 18).Value = 4.72
$ws.Cells.Item(65, 19).Value = "17/09/2023 09:02"
$ws.Cells.Item(65, 20).Value = 3.42
$ws.Cells.Item(65, 21).Value = "27/09/2023 18:51"
$ws.Cells.Item(65, 22).Value = "https://www.betexplorer.com/football/spain/laliga/villarreal-girona/80EuTg3A/"

# Row 66 (Indice=65)
$ws.Range("A62").Copy()
$ws.Range("A66").PasteSpecial($xlPasteFormats)
$ws.Range("E62").Copy()
$ws.Range("E66").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(66, 1).Value = 65
$ws.Cells.Item(66, 2).Value = "spain"
$ws.Cells.Item(66, 3).Value = "laliga"
$ws.Cells.Item(66, 4).Value = "2023-2024"
$ws.Cells.Item(66, 5).Value = 45196.89583333334
$ws.Cells.Item(66, 6).Value = "Cadiz CF"
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = "Rayo Vallecano"
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 2.59
$ws.Cells.Item(66, 11).Value = "17/09/2023 09:02"
$ws.Cells.Item(66, 12).Value = 2.72
$ws.Cells.Item(66, 13).Value = "27/09/2023 21:19"
$ws.Cells.Item(66, 14).Value = 3.06
$ws.Cells.Item(66, 15).Value = "17/09/2023 09:02"
$ws.Cells.Item(66, 16).Value = 3.11
$ws.Cells.Item(66, 17).Value = "27/09/2023 21:17"
$ws.Cells.Item(66, 18).Value = 3.11
$ws.Cells.Item(66, 19).Value = "17/09/2023 09:02"
$ws.Cells.Item(66, 20).Value = 3
$ws.Cells.Item(66, 21).Value = "27/09/2023 21:30"
$ws.Cells.Item(66, 22).Value = "https://www.betexplorer.com/football/spain/laliga/cadiz-rayo-vallecano/CEYt8hRp/"

# Row 67 (Indice=66)
$ws.Range("A62").Copy()
$ws.Range("A67").PasteSpecial($xlPasteFormats)
$ws.Range("E62").Copy()
$ws.Range("E67").PasteSpecial($xlPasteFormats)
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "spain"
$ws.Cells.Item(67, 3).Value = "laliga"
$ws.Cells.Item(67, 4).Value = "2023-2024"
$ws.Cells.Item(67, 5).Value = 45196.89583333334
$ws.Cells.Item(67, 6).Value = "Valencia"
$ws.Cells.Item(67, 7).Value = 0
$ws.Cells.Item(67, 8).Value = "Real Sociedad"
$ws.Cells.Item(67, 9).Value = 1
$ws.Cells.Item(67, 10).Value = 2.33
$ws.Cells.Item(67, 11).Value = "17/09/2023 09:02"
$ws.Cells.Item(67, 12).Value = 2.57
$ws.Cells.Item(67, 13).Value = "27/09/2023 21:27"
$ws.Cells.Item(67, 14).Value = 3.14
$ws.Cells.Item(67, 15).Value = "17/09/2023 09:02"
$ws.Cells.Item(67, 16).Value = 3.02
$ws.Cells.Item(67, 17).Value = "27/09/2023 21:27"
$ws.Cells.Item(67, 18).Value = 3.26
$ws.Cells.Item(67, 19).Value = "17/09/2023 09:02"
$ws.Cells.Item(67, 20).Value = 3.31
$ws.Cells.Item(67, 21).Value = "27/09/2023 21:27"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/spain/laliga/valencia-real-sociedad/M3IqSDIG/"

$excel.CutCopyMode = 0